# Update the "Generate Date"/"Datetime" timestamps recorded during report
# generation for the handback status report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the second data row.
$overview.Range("G3").Value = "2016-08-17 06:40:42"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the second data row.
$zhcn.Range("H3").Value = "2016-08-17 06:40:37"
$zhcn.Range("K3").Value = "2016-08-17 06:40:55"

# de-de sheet: same two columns for the second data row. The "Correspond
# Handoff Datetime" value mirrors the Overview sheet's generate date.
$dede.Range("H3").Value = "2016-08-17 06:40:42"
$dede.Range("K3").Value = "2016-08-17 06:41:06"
